$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding decimal-looking values that must remain plain text
# (force Text number format so Excel does not coerce them to floats)
$textCells = @("D5","D6","D9","D10","D12","D14","D15","D17","D19","D21","D22","D26","D27","D28","D30","D32","D34","D38","D39","D41","D43","D44","D45","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = '248.58'
$ws.Range("D6").Value = '0.666'
$ws.Range("D9").Value = '0.386'
$ws.Range("D10").Value = '0.0788'
$ws.Range("D12").Value = '15.83'
$ws.Range("D14").Value = '0.835'
$ws.Range("D15").Value = '5.80'
$ws.Range("D17").Value = '18.15'
$ws.Range("D19").Value = '75.17'
$ws.Range("D21").Value = '5.42'
$ws.Range("D22").Value = '238.15'
$ws.Range("D26").Value = '169.46'
$ws.Range("D27").Value = '9.40'
$ws.Range("D28").Value = '20.13'
$ws.Range("D30").Value = '4.84'
$ws.Range("D32").Value = '0.0622'
$ws.Range("D34").Value = '0.0908'
$ws.Range("D38").Value = '1.34'
$ws.Range("D39").Value = '0.106'
$ws.Range("D41").Value = '5.09'
$ws.Range("D43").Value = '17.41'
$ws.Range("D44").Value = '1.16'
$ws.Range("D45").Value = '96.64'
$ws.Range("D48").Value = '2.91'
$ws.Range("D49").Value = '6.88'
$ws.Range("D51").Value = '3.65'

# Remaining updated cells (already unambiguous text)
$ws.Range("D2").Value = '37.154.15'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '2.053.88'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("E7").Value = '  -6.50%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("E10").Value = '  -2.53%  '
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").Value = '2.355.24'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("E15").Value = '  +6.73%  '
$ws.Range("D16").Value = '2.054.16'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("E17").Value = '  +19.77%  '
$ws.Range("D18").Value = '37.206.01'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").Value = '0.0₃0902'
$ws.Range("E20").Value = '  -2.78%  '
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  +2.49%  '
$ws.Range("E25").Value = '  +5.39%  '
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("E27").Value = '  +1.22%  '
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("E31").Value = '  +1.96%  '
$ws.Range("E32").Value = '  -2.61%  '
$ws.Range("E33").Value = '  +2.75%  '
$ws.Range("E34").Value = '  +1.93%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -1.41%  '
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E38").Value = '  -0.56%  '
$ws.Range("E39").Value = '  -1.91%  '
$ws.Range("E40").Value = '  +12.34%  '
$ws.Range("E41").Value = '  +11.51%  '
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("E43").Value = '  -7.73%  '
$ws.Range("E44").Value = '  -0.81%  '
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("E46").Value = '  -3.45%  '
$ws.Range("D47").Value = '1.287.80'
$ws.Range("E47").Value = '  -1.81%  '
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("E49").Value = '  -0.80%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.248.38'
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("B51").Value = 'FTXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("E51").Value = '  -16.78%  '
